$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateUser")
$ws.Range("B2").Value = "Sat Nov 15 17:53:15 EST 2025"
$ws.Range("B3").Value = "Sat Nov 15 17:53:53 EST 2025"
$ws.Range("B4").Value = "Sat Nov 15 17:54:29 EST 2025"
$ws = $wb.Worksheets.Item("CreateUserSpChar")
$ws.Range("B2").Value = "Sat Nov 15 17:55:04 EST 2025"
$ws.Range("B3").Value = "Sat Nov 15 17:55:40 EST 2025"
$ws.Range("B4").Value = "Sat Nov 15 17:56:14 EST 2025"
$ws = $wb.Worksheets.Item("CreateUserSpCharError")
$ws.Range("B2").Value = "Sat Nov 15 17:56:54 EST 2025"
$ws.Range("B3").Value = "Sat Nov 15 17:57:11 EST 2025"
$ws.Range("B4").Value = "Sat Nov 15 17:57:28 EST 2025"
$ws.Range("B5").Value = "Sat Nov 15 17:57:44 EST 2025"
$ws.Range("B6").Value = "Sat Nov 15 17:58:02 EST 2025"
$ws.Range("B7").Value = "Sat Nov 15 17:58:19 EST 2025"
$ws = $wb.Worksheets.Item("FindUser")
$ws.Range("B2").Value = "Sat Nov 15 17:58:37 EST 2025"
$ws.Range("B3").Value = "Sat Nov 15 17:58:54 EST 2025"
$ws.Range("B4").Value = "Sat Nov 15 17:59:11 EST 2025"
$ws.Range("B5").Value = "Sat Nov 15 17:59:28 EST 2025"
$ws.Range("B6").Value = "Sat Nov 15 17:59:44 EST 2025"
$ws.Range("B7").Value = "Sat Nov 15 18:00:01 EST 2025"
$ws.Range("B8").Value = "Sat Nov 15 18:00:17 EST 2025"
$ws.Range("B9").Value = "Sat Nov 15 18:00:34 EST 2025"
$ws.Range("B10").Value = "Sat Nov 15 18:00:51 EST 2025"
$ws.Range("B11").Value = "Sat Nov 15 18:01:09 EST 2025"
$ws.Range("B12").Value = "Sat Nov 15 18:01:26 EST 2025"
$ws.Range("B13").Value = "Sat Nov 15 18:01:42 EST 2025"
$ws = $wb.Worksheets.Item("PassCase")
$ws.Range("B2").Value = "Sat Nov 15 18:02:04 EST 2025"
$ws = $wb.Worksheets.Item("UsernameCase")
$ws.Range("B2").Value = "Sat Nov 15 18:02:18 EST 2025"
$ws.Range("B3").Value = "Sat Nov 15 18:02:28 EST 2025"
$ws.Range("B4").Value = "Sat Nov 15 18:02:38 EST 2025"
$ws = $wb.Worksheets.Item("CreateUserPasswordSpChar")
$ws.Range("B2").Value = "Sat Nov 15 18:02:49 EST 2025"
$ws.Range("B3").Value = "Sat Nov 15 18:03:26 EST 2025"
$ws.Range("B4").Value = "Sat Nov 15 18:04:02 EST 2025"
$ws.Range("B5").Value = "Sat Nov 15 18:04:40 EST 2025"
$ws.Range("B6").Value = "Sat Nov 15 18:05:17 EST 2025"
$ws.Range("B7").Value = "Sat Nov 15 18:05:54 EST 2025"
$ws.Range("B8").Value = "Sat Nov 15 18:06:32 EST 2025"
$ws.Range("B9").Value = "Sat Nov 15 18:07:10 EST 2025"
$ws.Range("B10").Value = "Sat Nov 15 18:07:47 EST 2025"
$ws.Range("B11").Value = "Sat Nov 15 18:08:24 EST 2025"
$ws.Range("B12").Value = "Sat Nov 15 18:09:02 EST 2025"
$ws.Range("B13").Value = "Sat Nov 15 18:09:40 EST 2025"
$ws.Range("B14").Value = "Sat Nov 15 18:10:17 EST 2025"
$ws.Range("B15").Value = "Sat Nov 15 18:10:54 EST 2025"
$ws.Range("B16").Value = "Sat Nov 15 18:11:31 EST 2025"
$ws.Range("B17").Value = "Sat Nov 15 18:12:09 EST 2025"
$ws = $wb.Worksheets.Item("ModifyUserPwd")
$ws.Range("B2").Value = "Sat Nov 15 18:14:06 EST 2025"
$ws.Range("B3").Value = "Sat Nov 15 18:14:34 EST 2025"
$ws.Range("B4").Value = "Sat Nov 15 18:15:01 EST 2025"
$ws.Range("B5").Value = "Sat Nov 15 18:15:27 EST 2025"
$ws.Range("B6").Value = "Sat Nov 15 18:15:54 EST 2025"
$ws.Range("B7").Value = "Sat Nov 15 18:16:21 EST 2025"
$ws.Range("B8").Value = "Sat Nov 15 18:16:45 EST 2025"
$ws = $wb.Worksheets.Item("ModifyUser")
$ws.Range("B2").Value = "Sat Nov 15 18:17:10 EST 2025"
$ws.Range("B3").Value = "Sat Nov 15 18:17:58 EST 2025"
$ws = $wb.Worksheets.Item("CreateUserSCFNameErr")
$ws.Range("B2").Value = "Sat Nov 15 18:18:47 EST 2025"
$ws.Range("B3").Value = "Sat Nov 15 18:19:09 EST 2025"
$ws.Range("B4").Value = "Sat Nov 15 18:19:30 EST 2025"
$ws.Range("B5").Value = "Sat Nov 15 18:19:52 EST 2025"
$ws = $wb.Worksheets.Item("CreateUserSCLNameErr")
$ws.Range("B2").Value = "Sat Nov 15 18:20:13 EST 2025"
$ws.Range("B3").Value = "Sat Nov 15 18:20:35 EST 2025"
$ws.Range("B4").Value = "Sat Nov 15 18:20:56 EST 2025"
$ws.Range("B5").Value = "Sat Nov 15 18:21:18 EST 2025"
$ws = $wb.Worksheets.Item("CreateUserErrors")
$ws.Range("B2").Value = "Sat Nov 15 18:21:40 EST 2025"
$ws.Range("B3").Value = "Sat Nov 15 18:22:01 EST 2025"
$ws.Range("B4").Value = "Sat Nov 15 18:22:22 EST 2025"
$ws.Range("B5").Value = "Sat Nov 15 18:22:43 EST 2025"
$ws.Range("B6").Value = "Sat Nov 15 18:23:05 EST 2025"
$ws.Range("B8").Value = "Sat Nov 15 18:23:26 EST 2025"
$ws.Range("B9").Value = "Sat Nov 15 18:23:47 EST 2025"
$ws.Range("B10").Value = "Sat Nov 15 18:24:09 EST 2025"
$ws.Range("B11").Value = "Sat Nov 15 18:24:30 EST 2025"
$ws.Range("B12").Value = "Sat Nov 15 18:24:52 EST 2025"
$ws = $wb.Worksheets.Item("AddDeleteRole")
$ws.Range("B2").Value = "Sat Nov 15 18:26:12 EST 2025"
$ws.Range("B3").Value = "Sat Nov 15 18:26:59 EST 2025"
$ws.Range("B4").Value = "Sat Nov 15 18:27:48 EST 2025"
$ws.Range("B5").Value = "Sat Nov 15 18:28:36 EST 2025"
$ws = $wb.Worksheets.Item("SearchRole")
$ws.Range("B2").Value = "Sat Nov 15 18:29:24 EST 2025"
$ws.Range("B3").Value = "Sat Nov 15 18:29:38 EST 2025"
$ws.Range("B4").Value = "Sat Nov 15 18:29:52 EST 2025"
$ws.Range("B5").Value = "Sat Nov 15 18:30:05 EST 2025"
